$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 is a text code with a leading zero ("001") - force text format so it
# is not coerced into the number 1, then restore the cell style so no
# stray number-format style is left attached to the cell.
$jCell = $ws.Cells.Item(2, 10)
$jCell.NumberFormat = "@"
$jCell.Value = "001"
$jCell.Style = "Normal"

# N2 is a plain text timestamp string (inline string in the source data).
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# Numeric metrics for row 2.
$ws.Range("O2").Value = -14067931.59
$ws.Range("P2").Value = -32.655542313
$ws.Range("Q2").Value = 145268344.38
$ws.Range("R2").Value = 337.2078216537
$ws.Range("S2").Value = 47858180.31
$ws.Range("T2").Value = 111.0920125064
$ws.Range("U2").Value = 2895869.73
$ws.Range("V2").Value = 6.7221109156
$ws.Range("Y2").Value = 12149132.79
$ws.Range("Z2").Value = 28.2014820268
$ws.Range("AA2").Value = -30998727.34
$ws.Range("AB2").Value = -71.95658052669999
$ws.Range("AC2").Value = -43079767.15
